$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: Coin (B), Link (C), Price (D), Volume 1h (E), Hora (G)
# Price/Volume/Hora are stored as text in the source sheet (leading zeros /
# "%" suffix / literal "--" must survive), so force Text format before writing.
$rows = @(
    @{ Row=2; D="306.92"; E="0.17%"; G="1" }
    @{ Row=3; D="41.10"; E="2.56%"; G="1" }
    @{ Row=4; D="5.122"; E="2.13%"; G="1" }
    @{ Row=5; D="0.07611"; E="-0.89%"; G="1" }
    @{ Row=6; B="FTXToken"; C="https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D="1.631"; E="0.85%"; G="1" }
    @{ Row=7; B="BTSEToken"; C="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D="2.533"; E="-0.22%"; G="1" }
    @{ Row=8; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="0.9050"; E="-0.04%"; G="1" }
    @{ Row=9; B="LiechtensteinCryptoassetsExchange"; C="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D="0.1064"; E="7.66%"; G="1" }
    @{ Row=10; B="WazirX"; C="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D="0.1752"; E="1.68%"; G="1" }
    @{ Row=11; B="MandalaExchangeToken"; C="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D="0.09156"; E="1.41%"; G="1" }
    @{ Row=12; B="BitrueCoin"; C="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D="0.04196"; E="-5.01%"; G="1" }
    @{ Row=13; B="BitMartToken"; C="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D="0.1051"; E="-0.60%"; G="1" }
    @{ Row=14; B="BitForexToken"; C="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D="0.001255"; E="-0.94%"; G="1" }
    @{ Row=15; B="CoinExToken"; C="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D="0.04174"; E="-0.36%"; G="1" }
    @{ Row=16; D="0.005834"; E="0.73%"; G="1" }
    @{ Row=17; D="3.352"; E="-0.19%"; G="1" }
    @{ Row=18; B="GateToken"; C="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D="4.252"; E="-0.16%"; G="1" }
    @{ Row=19; B="BitpandaEcosystemToken"; C="https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D="0.3275"; E="-2.68%"; G="1" }
    @{ Row=20; B="MCDex"; C="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; D="6.574"; E="-6.87%"; G="1" }
    @{ Row=21; B="ProBitToken"; C="https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; D="0.1360"; E="1.45%"; G="1" }
    @{ Row=22; B="ZBToken"; C="https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; D="0.2721"; E="-5.11%"; G="1" }
    @{ Row=23; D="0.001223"; E="2.61%"; G="1" }
    @{ Row=24; D="0.004070"; E="-0.27%"; G="1" }
    @{ Row=25; D="0.0001299"; E="6.40%"; G="1" }
    @{ Row=26; D="0.0003005"; E="0.80%"; G="1" }
    @{ Row=27; G="1" }
    @{ Row=28; G="1" }
    @{ Row=29; G="1" }
    @{ Row=30; G="1" }
    @{ Row=31; G="1" }
    @{ Row=32; G="1" }
    @{ Row=33; G="1" }
    @{ Row=34; G="1" }
    @{ Row=35; G="1" }
    @{ Row=36; G="1" }
    @{ Row=37; G="1" }
    @{ Row=38; D="0.02371"; E="1.74%"; G="1" }
    @{ Row=39; D="0.05168"; E="0.75%"; G="1" }
    @{ Row=40; D="0.007764"; E="-1.75%"; G="1" }
    @{ Row=41; E="-2.32%"; G="1" }
    @{ Row=42; D="0.006960"; E="0.46%"; G="1" }
    @{ Row=43; D="0.001917"; E="1.10%"; G="1" }
    @{ Row=44; D="0.008578"; E="7.30%"; G="1" }
    @{ Row=45; D="0.3045"; E="-8.06%"; G="1" }
    @{ Row=46; D="0.00006376"; E="-4.38%"; G="1" }
    @{ Row=47; D="0.00000000749"; E="-0.15%"; G="1" }
    @{ Row=48; B="CoinbaseStockToken"; C="https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"; D="0.004397"; E="6.86%"; G="1" }
    @{ Row=49; B="BOLO"; C="https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"; D="0.009547"; E="338.25%"; G="1" }
    @{ Row=50; D="0.00002098"; E="-0.15%"; G="1" }
    @{ Row=51; D="0.0001998"; E="-0.15%"; G="1" }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($r.ContainsKey("B")) { $ws.Range("B$row").Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Range("C$row").Value = $r.C }
    if ($r.ContainsKey("D")) {
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $r.D
    }
    if ($r.ContainsKey("E")) {
        $ws.Range("E$row").NumberFormat = "@"
        $ws.Range("E$row").Value = $r.E
    }
    # Hora switches from "0" to "1" for every data row in this update
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = "1"
}
